$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1282.25
$ws.Range("J4").Value = 150
$ws.Range("L4").Value = 150
$ws.Range("N4").Value = -378
$ws.Range("H19").Value = 485.125
$ws.Range("I19").Value = 236.25
$ws.Range("J19").Value = 734
$ws.Range("K19").Value = 236.25
$ws.Range("L19").Value = 734
$ws.Range("M19").Value = -61.25
$ws.Range("N19").Value = -1084
$ws.Range("H40").Value = 2106.3684
$ws.Range("I40").Value = 2070.8462
$ws.Range("J40").Value = 2183.3333
$ws.Range("K40").Value = 2070.8462
$ws.Range("L40").Value = 2183.3333
$ws.Range("M40").Value = -1895.8462
$ws.Range("N40").Value = -2533.3333
$ws.Range("H51").Value = 5640
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 6550
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 6550
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -7518
$ws.Range("H55").Value = 218.5
$ws.Range("I55").Value = 210
$ws.Range("J55").Value = 261
$ws.Range("K55").Value = 210
$ws.Range("L55").Value = 261
$ws.Range("M55").Value = 4
$ws.Range("N55").Value = -689
$ws.Range("H113").Value = 103850.4
$ws.Range("I113").Value = 128688
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 128688
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -125434
$ws.Range("N113").Value = -11008
$ws.Range("H141").Value = 2182.0952
$ws.Range("J141").Value = 6492.5454
$ws.Range("L141").Value = 19477.6362
$ws.Range("N141").Value = -29837.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 66055.69
$ws.Range("I2").Value = 115176.78
$ws.Range("K2").Value = 115176.78
$ws.Range("M2").Value = -115063.78
$ws.Range("H61").Value = 999.3492
$ws.Range("I61").Value = 956.15
$ws.Range("J61").Value = 1863.3334
$ws.Range("K61").Value = 956.15
$ws.Range("L61").Value = 1863.3334
$ws.Range("M61").Value = -744.15
$ws.Range("N61").Value = -2287.3334
$ws.Range("H63").Value = 10320.714
$ws.Range("I63").Value = 19141.428
$ws.Range("J63").Value = 1500
$ws.Range("K63").Value = 19141.428
$ws.Range("L63").Value = 1500
$ws.Range("M63").Value = -18455.428
$ws.Range("N63").Value = -2872
$ws.Range("H66").Value = 10320.714
$ws.Range("I66").Value = 19141.428
$ws.Range("J66").Value = 1500
$ws.Range("K66").Value = 95707.14
$ws.Range("L66").Value = 7500
$ws.Range("M66").Value = -92275.14
$ws.Range("N66").Value = -14364
$ws.Range("H116").Value = 66055.69
$ws.Range("I116").Value = 115176.78
$ws.Range("K116").Value = 115176.78
$ws.Range("M116").Value = -112882.78
$ws.Range("H122").Value = 7995.1055
$ws.Range("I122").Value = 9367.134
$ws.Range("J122").Value = 2850
$ws.Range("K122").Value = 28101.402
$ws.Range("L122").Value = 8550
$ws.Range("M122").Value = -25651.402
$ws.Range("N122").Value = -13450
$ws.Range("H132").Value = 1571.6666
$ws.Range("I132").Value = 1322.9836
$ws.Range("J132").Value = 3467.875
$ws.Range("K132").Value = 3968.9508
$ws.Range("L132").Value = 10403.625
$ws.Range("M132").Value = -1438.9508
$ws.Range("N132").Value = -15463.625
$ws.Range("H136").Value = 999.3492
$ws.Range("I136").Value = 956.15
$ws.Range("J136").Value = 1863.3334
$ws.Range("K136").Value = 2868.45
$ws.Range("L136").Value = 5590.0002
$ws.Range("M136").Value = -318.4499999999998
$ws.Range("N136").Value = -10690.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 66055.69
$ws.Range("I3").Value = 115176.78
$ws.Range("K3").Value = 115176.78
$ws.Range("M3").Value = -115062.78
$ws.Range("H20").Value = 3999.8
$ws.Range("I20").Value = 3999.75
$ws.Range("K20").Value = 3999.75
$ws.Range("M20").Value = -3752.75
$ws.Range("H134").Value = 14927704
$ws.Range("I134").Value = 19232508
$ws.Range("J134").Value = 4381.067
$ws.Range("K134").Value = 57697524
$ws.Range("L134").Value = 13143.201
$ws.Range("M134").Value = -57694989
$ws.Range("N134").Value = -18213.201

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 986.44446
$ws.Range("I16").Value = 449.66666
$ws.Range("K16").Value = 449.66666
$ws.Range("M16").Value = -162.66666
$ws.Range("H31").Value = 1875.1951
$ws.Range("I31").Value = 1060.2812
$ws.Range("K31").Value = 1060.2812
$ws.Range("M31").Value = -765.2811999999999
$ws.Range("H34").Value = 1875.1951
$ws.Range("I34").Value = 1060.2812
$ws.Range("K34").Value = 1060.2812
$ws.Range("M34").Value = -858.2811999999999
$ws.Range("H99").Value = 62500000
$ws.Range("I99").Value = 62500000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 62500000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -62498502
$ws.Range("N99").ClearContents()
$ws.Range("H112").Value = 25490
$ws.Range("J112").Value = 25490
$ws.Range("L112").Value = 25490
$ws.Range("N112").Value = -28444
$ws.Range("H113").Value = 986.44446
$ws.Range("I113").Value = 449.66666
$ws.Range("K113").Value = 449.66666
$ws.Range("M113").Value = 1720.33334
$ws.Range("H126").Value = 62500000
$ws.Range("I126").Value = 62500000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 187500000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -187497530
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1396.1167
$ws.Range("I132").Value = 1244.3922
$ws.Range("J132").Value = 2255.889
$ws.Range("K132").Value = 3733.1766
$ws.Range("L132").Value = 6767.667
$ws.Range("M132").Value = -1203.1766
$ws.Range("N132").Value = -11827.667
$ws.Range("H134").Value = 1622.3735
$ws.Range("I134").Value = 1026.2787
$ws.Range("J134").Value = 3275.182
$ws.Range("K134").Value = 3078.8361
$ws.Range("L134").Value = 9825.545999999998
$ws.Range("M134").Value = -543.8361000000004
$ws.Range("N134").Value = -14895.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6378
$ws.Range("I3").Value = 2331.6667
$ws.Range("J3").Value = 18517
$ws.Range("K3").Value = 6995.000100000001
$ws.Range("L3").Value = 55551
$ws.Range("M3").Value = -6883.000100000001
$ws.Range("N3").Value = -55775
$ws.Range("H16").Value = 911.5714
$ws.Range("I16").Value = 927
$ws.Range("K16").Value = 2781
$ws.Range("M16").Value = -2608
$ws.Range("H40").Value = 520
$ws.Range("I40").Value = 350
$ws.Range("J40").Value = 775
$ws.Range("K40").Value = 1400
$ws.Range("L40").Value = 3100
$ws.Range("M40").Value = -1331
$ws.Range("N40").Value = -3238
$ws.Range("H131").Value = 2224.3333
$ws.Range("J131").Value = 2340.0967
$ws.Range("L131").Value = 7020.2901
$ws.Range("N131").Value = -17100.2901

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6925.8184
$ws.Range("I70").Value = 7684.9165
$ws.Range("J70").Value = 6014.9
$ws.Range("K70").Value = 7684.9165
$ws.Range("L70").Value = 6014.9
$ws.Range("M70").Value = -7414.9165
$ws.Range("N70").Value = -6554.9
$ws.Range("H73").Value = 6925.8184
$ws.Range("I73").Value = 7684.9165
$ws.Range("J73").Value = 6014.9
$ws.Range("K73").Value = 7684.9165
$ws.Range("L73").Value = 6014.9
$ws.Range("M73").Value = -6748.9165
$ws.Range("N73").Value = -7886.9
$ws.Range("H122").Value = 585952.9
$ws.Range("I122").Value = 654747.4
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 1964242.2
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -1961792.2
$ws.Range("N122").Value = -8500
$ws.Range("H132").Value = 2187.4626
$ws.Range("I132").Value = 1908.3392
$ws.Range("J132").Value = 3608.4546
$ws.Range("K132").Value = 5725.017599999999
$ws.Range("L132").Value = 10825.3638
$ws.Range("M132").Value = -3195.017599999999
$ws.Range("N132").Value = -15885.3638
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4763347.5
$ws.Range("J16").Value = 734.5
$ws.Range("L16").Value = 734.5
$ws.Range("N16").Value = -1074.5
$ws.Range("H82").Value = 1062.875
$ws.Range("I82").Value = 971.4286
$ws.Range("J82").Value = 1703
$ws.Range("K82").Value = 971.4286
$ws.Range("L82").Value = 1703
$ws.Range("M82").Value = -610.4286
$ws.Range("N82").Value = -2425
$ws.Range("H85").Value = 1062.875
$ws.Range("I85").Value = 971.4286
$ws.Range("J85").Value = 1703
$ws.Range("K85").Value = 971.4286
$ws.Range("L85").Value = 1703
$ws.Range("M85").Value = 276.5714
$ws.Range("N85").Value = -4199
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H136").Value = 2974.7747
$ws.Range("I136").Value = 1932.3167
$ws.Range("K136").Value = 5796.9501
$ws.Range("M136").Value = -3246.9501
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1002307
$ws.Range("I81").Value = 1430517.2
$ws.Range("J81").Value = 3149.8333
$ws.Range("K81").Value = 2861034.4
$ws.Range("L81").Value = 6299.6666
$ws.Range("M81").Value = -2859973.4
$ws.Range("N81").Value = -8421.6666
$ws.Range("H84").Value = 1002307
$ws.Range("I84").Value = 1430517.2
$ws.Range("J84").Value = 3149.8333
$ws.Range("K84").Value = 14305172
$ws.Range("L84").Value = 31498.333
$ws.Range("M84").Value = -14299868
$ws.Range("N84").Value = -42106.333
$ws.Range("H100").Value = 1315.9231
$ws.Range("I100").Value = 1638
$ws.Range("J100").Value = 800.6
$ws.Range("K100").Value = 3276
$ws.Range("L100").Value = 1601.2
$ws.Range("M100").Value = -2735
$ws.Range("N100").Value = -2683.2
$ws.Range("H107").Value = 2646473.2
$ws.Range("I107").Value = 4630454
$ws.Range("J107").Value = 1165.5555
$ws.Range("K107").Value = 13891362
$ws.Range("L107").Value = 3496.6665
$ws.Range("M107").Value = -13889442
$ws.Range("N107").Value = -7336.666499999999
$ws.Range("H132").Value = 7577502
$ws.Range("I132").Value = 9435505
$ws.Range("J132").Value = 2564.1538
$ws.Range("K132").Value = 28306515
$ws.Range("L132").Value = 7692.4614
$ws.Range("M132").Value = -28303985
$ws.Range("N132").Value = -12752.4614
$ws.Range("H136").Value = 14848.194
$ws.Range("I136").Value = 15996.303
$ws.Range("K136").Value = 47988.909
$ws.Range("M136").Value = -45438.909
